$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for Price cells that would otherwise be parsed as numbers,
# matching the workbook's inline-string (text) cell type.
$numericPriceCells = @("D4","D5","D6","D7","D9","D10","D11","D13","D14","D15","D16","D17","D19","D20","D21","D22","D23","D24","D26","D27","D28","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $numericPriceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated values
$ws.Range("D2").Value = "20.559.66"
$ws.Range("E2").Value = "  +1.67%  "
$ws.Range("D3").Value = "1.472.77"
$ws.Range("E3").Value = "  +2.12%  "
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "0.9587"
$ws.Range("E5").Value = "  +5.21%  "
$ws.Range("D6").Value = "277.20"
$ws.Range("E6").Value = "  -0.24%  "
$ws.Range("D7").Value = "0.3618"
$ws.Range("E7").Value = "  -0.83%  "
$ws.Range("E8").Value = "  -1.13%  "
$ws.Range("D9").Value = "39.69"
$ws.Range("E9").Value = "  +1.97%  "
$ws.Range("D10").Value = "1.074"
$ws.Range("E10").Value = "  +5.33%  "
$ws.Range("D11").Value = "0.06660"
$ws.Range("E11").Value = "  +2.15%  "
$ws.Range("E12").Value = "  +0.04%  "
$ws.Range("D13").Value = "5.521"
$ws.Range("E13").Value = "  +2.47%  "
$ws.Range("D14").Value = "18.19"
$ws.Range("E14").Value = "  +3.33%  "
$ws.Range("D15").Value = "6.175"
$ws.Range("E15").Value = "  +2.09%  "
$ws.Range("D16").Value = "0.9587"
$ws.Range("E16").Value = "  +2.02%  "
$ws.Range("D17").Value = "0.00001027"
$ws.Range("E17").Value = "  +1.01%  "
$ws.Range("D18").Value = "1.472.96"
$ws.Range("E18").Value = "  +2.12%  "
$ws.Range("D19").Value = "0.05940"
$ws.Range("E19").Value = "  +5.44%  "
$ws.Range("D20").Value = "68.87"
$ws.Range("E20").Value = "  +1.26%  "
$ws.Range("D21").Value = "5.498"
$ws.Range("E21").Value = "  +2.18%  "
$ws.Range("D22").Value = "14.55"
$ws.Range("E22").Value = "  +1.07%  "
$ws.Range("D23").Value = "11.17"
$ws.Range("E23").Value = "  +3.45%  "
$ws.Range("D24").Value = "2.266"
$ws.Range("E24").Value = "  +0.25%  "
$ws.Range("D25").Value = "20.552.47"
$ws.Range("E25").Value = "  +1.46%  "
$ws.Range("D26").Value = "143.01"
$ws.Range("E26").Value = "  +3.70%  "
$ws.Range("D27").Value = "2.129"
$ws.Range("E27").Value = "  -1.56%  "
$ws.Range("D28").Value = "17.15"
$ws.Range("E28").Value = "  +1.40%  "
$ws.Range("D29").Value = "1.632.95"
$ws.Range("E29").Value = "  +2.33%  "
$ws.Range("D30").Value = "113.97"
$ws.Range("E30").Value = "  +3.73%  "
$ws.Range("D31").Value = "3.896"
$ws.Range("E31").Value = "  +1.66%  "
$ws.Range("D32").Value = "4.961"
$ws.Range("E32").Value = "  +2.86%  "
$ws.Range("D33").Value = "0.08019"
$ws.Range("D34").Value = "0.8050"
$ws.Range("E34").Value = "  +0.08%  "
$ws.Range("D35").Value = "1.517"
$ws.Range("E35").Value = "  +4.90%  "
$ws.Range("D36").Value = "1.214"
$ws.Range("E36").Value = "  +6.92%  "
$ws.Range("D37").Value = "0.05765"
$ws.Range("E37").Value = "  -2.52%  "
$ws.Range("D38").Value = "4.743"
$ws.Range("E38").Value = "  +1.63%  "
$ws.Range("D39").Value = "0.02057"
$ws.Range("E39").Value = "  +3.33%  "
$ws.Range("B40").Value = "Aptos"
$ws.Range("C40").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D40").Value = "10.43"
$ws.Range("E40").Value = "  +2.84%  "
$ws.Range("B41").Value = "Frax"
$ws.Range("C41").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D41").Value = "0.9593"
$ws.Range("E41").Value = "  +3.17%  "
$ws.Range("D42").Value = "0.1873"
$ws.Range("E42").Value = "  +1.87%  "
$ws.Range("D43").Value = "7.405"
$ws.Range("E43").Value = "  +3.47%  "
$ws.Range("D44").Value = "0.5286"
$ws.Range("E44").Value = "  +1.17%  "
$ws.Range("D45").Value = "3.517"
$ws.Range("E45").Value = "  -0.08%  "
$ws.Range("D46").Value = "12.19"
$ws.Range("E46").Value = "  +1.13%  "
$ws.Range("D47").Value = "118.60"
$ws.Range("E47").Value = "  -0.40%  "
$ws.Range("D48").Value = "0.5207"
$ws.Range("E48").Value = "  +1.51%  "
$ws.Range("D49").Value = "1.820"
$ws.Range("E49").Value = "  +3.70%  "
$ws.Range("D50").Value = "0.06475"
$ws.Range("E50").Value = "  +2.10%  "
$ws.Range("D51").Value = "0.9879"
$ws.Range("E51").Value = "  -0.59%  "

# Restore default styling so the text-forcing NumberFormat change above
# does not leave a residual style index on these cells.
foreach ($addr in $numericPriceCells) {
    $ws.Range($addr).Style = "Normal"
}
